$d = $word.ActiveDocument

# --- Locate the pieces of text we need to touch -----------------------
# "Version" currently spans two runs ("Versi" + "on"); find its extent.
$rVersion = $d.Content
$rVersion.Find.Execute("Version", $false, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)
$versionStart = $rVersion.Start
$versionEnd = $rVersion.End

# " 2" is the run right after "Version" (before the _GoBack bookmark).
$rNum = $d.Content
$rNum.Find.Execute(" 2", $false, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0)
$numStart = $rNum.Start
$numEnd = $rNum.End

# --- 1) Merge "Versi" + "on" runs into a single "Version" run ---------
# A direct same-value assignment is a no-op for the engine (it already
# reads "Version"), so nudge it through a temporary value first to force
# the runs to merge, then set the final text.
$rTmp1 = $d.Range($versionStart, $versionEnd)
$rTmp1.Text = "VersionZZ"
$rFinal1 = $d.Range($versionStart, $versionStart + 9)
$rFinal1.Text = "Version"

# --- 2) Change " 2" into " 1." (absorbing the trailing period) --------
$rTmp2 = $d.Range($numStart, $numEnd)
$rTmp2.Text = " X"
$rFinal2 = $d.Range($numStart, $numStart + 2)
$rFinal2.Text = " 1."

# --- 3) Remove the now-redundant trailing "." run after the bookmark --
$periodStart = $numStart + 3
$rPeriod = $d.Range($periodStart, $periodStart + 1)
$rPeriod.Delete()
